# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "SANDRA MILENA CASTILLO ORTIZ" (CC 1143343913) rows move up to the top
# of the data block (rows 16-18), listed in descending period order
# (2212, 2211, 2210), and her last period's values are refreshed
# (F=40000, G=1300000 for all three rows now instead of 9333/40000/40000).
# The "ADA LUZ UTRIA NAVARRO" (CC 1048943696) row drops to the bottom
# (row 19), keeping her original figures (period 1808, F=16666, G=1000000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: SANDRA MILENA CASTILLO ORTIZ - period 2212
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143343913"
$ws.Range("D16").Value = "SANDRA MILENA CASTILLO ORTIZ"
$ws.Range("E16").Value = "2212"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1300000

# Row 17: SANDRA MILENA CASTILLO ORTIZ - period 2211
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143343913"
$ws.Range("D17").Value = "SANDRA MILENA CASTILLO ORTIZ"
$ws.Range("E17").Value = "2211"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1300000

# Row 18: SANDRA MILENA CASTILLO ORTIZ - period 2210
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143343913"
$ws.Range("D18").Value = "SANDRA MILENA CASTILLO ORTIZ"
$ws.Range("E18").Value = "2210"
$ws.Range("F18").Value = 9333
$ws.Range("G18").Value = 1300000

# Row 19: ADA LUZ UTRIA NAVARRO - period 1808
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1048943696"
$ws.Range("D19").Value = "ADA LUZ UTRIA NAVARRO"
$ws.Range("E19").Value = "1808"
$ws.Range("F19").Value = 16666
$ws.Range("G19").Value = 1000000
